# Experiment 2b participant_info.xlsx
# Commit: "update % DC in E2b"
#
# A new column "DC_lett" is inserted before the old column G (Perc_DC_D),
# shifting the old G/H/I (Perc_DC_D / Perc_DC_O / Perc_DC_M) to H/I/J.
# The new column is populated with 0/1 values per participant, a summary
# formula is added for it (and for column D), a comment's wording is
# tweaked, and the two comments that used to sit on the old G1/H1 headers
# move with their cells to H1/I1. The hidden _FilterDatabase name is
# widened to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture the text of the comments that will need to move, and the
#    corrected wording for the F1 comment, before we shuffle columns.
# ---------------------------------------------------------------------
$oldG1CommentText = $ws.Range("G1").Comment.Text()
$oldH1CommentText = $ws.Range("H1").Comment.Text()
$newF1CommentText = "awarenes of letter mask changes"

# ---------------------------------------------------------------------
# 2. Insert the new column at G, pushing the old G/H/I one place right.
# ---------------------------------------------------------------------
$ws.Columns("G:G").Insert()

# Header for the freshly inserted column.
$ws.Range("G1").Value = "DC_lett"

# ---------------------------------------------------------------------
# 3. Fill the new DC_lett column (rows 2-61) with the recorded values.
# ---------------------------------------------------------------------
$dcLett = @(1,0,0,1,0,0,0,0,1,0,0,0,1,0,0,0,1,0,1,1,1,0,0,1,1,0,1,1,1,1,0,0,1,1,1,0,1,1,1,1,1,1,1,0,1,1,0,0,0,0,0,0,1,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $dcLett.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $dcLett[$i]
}

# ---------------------------------------------------------------------
# 4. Summary formulas on row 63: proportion of "1" for D (DC_M) and the
#    newly inserted G (DC_lett) columns.
# ---------------------------------------------------------------------
$ws.Range("D63").Formula = "=SUM(D2:D61)/60"
$ws.Range("G63").Formula = "=SUM(G2:G61)/60"

# All of the DC_lett column (data + the new summary cell) is centred,
# matching the rest of the sheet's data columns.
$ws.Range("G1:G63").HorizontalAlignment = -4108  # xlCenter

# ---------------------------------------------------------------------
# 5. Fix up comments. Inserting a column does not relocate existing
#    cell comments in this runtime, so move them by hand:
#      F1 -> reworded text (same cell)
#      G1 (old wording, still physically on G1) -> moves to H1
#      H1 (old wording, still physically on H1) -> moves to I1
# ---------------------------------------------------------------------
$ws.Range("F1").Comment.Delete()
$ws.Range("F1").AddComment($newF1CommentText)

$ws.Range("G1").Comment.Delete()
$ws.Range("H1").Comment.Delete()

$ws.Range("H1").AddComment($oldG1CommentText)
$ws.Range("I1").AddComment($oldH1CommentText)

# ---------------------------------------------------------------------
# 6. Widen the hidden _FilterDatabase defined name to include the new
#    column (Excel normally does this automatically for a column
#    inserted inside the filtered range).
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$J`$66"
    }
}

# ---------------------------------------------------------------------
# 7. Restore the view state recorded in the edited workbook.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("F54").Select()
